# This script updates the "dSF" column (column F) values on Sheet1 to
# reflect a repull/recalculation of data for several rows, per commit
# message: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F ("dSF")
$updates = @{
    4  = -1
    5  = -2
    9  = 3
    10 = 4
    12 = -5
    13 = 4
    14 = -3
    15 = -3
    16 = -1
    17 = -1
    18 = 1
    19 = -1
    20 = -3
    21 = -4
    22 = 5
    23 = 1
    24 = 2
    25 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
